# Updated cryptos list values (Price / Volume(1h)) and the
# TheSandbox / WEMIXTOKEN row-order fix, per the Thu Mar 30 23:13:22 UTC 2023
# GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    # Column D prices are plain text (e.g. "20.50", "0.1060"); a bare
    # .Value assignment of a numeric-looking string auto-converts to a
    # Double (dropping trailing zeros / switching type). Force Text first,
    # then restore the default 'Normal' style so no stray number format
    # sticks around on the cell.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Row 2
$ws.Range("D2").Value = '27.978.51'
$ws.Range("E2").Value = '  -1.65%  '

# Row 3
$ws.Range("D3").Value = '1.789.99'
$ws.Range("E3").Value = '  -0.32%  '

# Row 4
$ws.Range("E4").Value = '  +0.26%  '

# Row 5
Set-TextValue $ws.Range("D5") '316.98'
$ws.Range("E5").Value = '  +0.86%  '

# Row 6
Set-TextValue $ws.Range("D6") '1.002'
$ws.Range("E6").Value = '  +0.19%  '

# Row 7
Set-TextValue $ws.Range("D7") '0.5344'
$ws.Range("E7").Value = '  -2.08%  '

# Row 8
Set-TextValue $ws.Range("D8") '0.3764'
$ws.Range("E8").Value = '  -1.48%  '

# Row 9
Set-TextValue $ws.Range("D9") '0.07409'
$ws.Range("E9").Value = '  -2.59%  '

# Row 10
Set-TextValue $ws.Range("D10") '41.76'
$ws.Range("E10").Value = '  -1.74%  '

# Row 11
Set-TextValue $ws.Range("D11") '1.089'
$ws.Range("E11").Value = '  -3.13%  '

# Row 12
Set-TextValue $ws.Range("D12") '1.003'
$ws.Range("E12").Value = '  +0.30%  '

# Row 13
Set-TextValue $ws.Range("D13") '20.50'
$ws.Range("E13").Value = '  -3.30%  '

# Row 14
Set-TextValue $ws.Range("D14") '6.103'
$ws.Range("E14").Value = '  -1.61%  '

# Row 15
$ws.Range("D15").Value = '1.785.25'
$ws.Range("E15").Value = '  -0.76%  '

# Row 16
Set-TextValue $ws.Range("D16") '7.216'
$ws.Range("E16").Value = '  -2.52%  '

# Row 17
Set-TextValue $ws.Range("D17") '88.92'
$ws.Range("E17").Value = '  -2.92%  '

# Row 18
Set-TextValue $ws.Range("D18") '0.00001056'
$ws.Range("E18").Value = '  -1.47%  '

# Row 19
Set-TextValue $ws.Range("D19") '0.06483'
$ws.Range("E19").Value = '  +0.36%  '

# Row 20
Set-TextValue $ws.Range("D20") '1.001'
$ws.Range("E20").Value = '  +0.13%  '

# Row 21
$ws.Range("E21").Value = '  -0.85%  '

# Row 22
Set-TextValue $ws.Range("D22") '5.887'
$ws.Range("E22").Value = '  -1.39%  '

# Row 23
$ws.Range("D23").Value = '28.017.21'
$ws.Range("E23").Value = '  -1.50%  '

# Row 24
Set-TextValue $ws.Range("D24") '11.14'
$ws.Range("E24").Value = '  -2.61%  '

# Row 25
Set-TextValue $ws.Range("D25") '2.095'
$ws.Range("E25").Value = '  -1.06%  '

# Row 26
Set-TextValue $ws.Range("D26") '155.76'
$ws.Range("E26").Value = '  -2.80%  '

# Row 27
Set-TextValue $ws.Range("D27") '20.26'
$ws.Range("E27").Value = '  -2.29%  '

# Row 28
$ws.Range("D28").Value = '1.991.98'
$ws.Range("E28").Value = '  -0.60%  '

# Row 29
Set-TextValue $ws.Range("D29") '2.278'
$ws.Range("E29").Value = '  -5.24%  '

# Row 30
Set-TextValue $ws.Range("D30") '120.86'
$ws.Range("E30").Value = '  -1.97%  '

# Row 31
Set-TextValue $ws.Range("D31") '1.111'
$ws.Range("E31").Value = '  -1.16%  '

# Row 32
Set-TextValue $ws.Range("D32") '0.1060'
$ws.Range("E32").Value = '  +3.16%  '

# Row 33
Set-TextValue $ws.Range("D33") '3.658'
$ws.Range("E33").Value = '  -0.45%  '

# Row 34
$ws.Range("E34").Value = '  -4.00%  '

# Row 35
Set-TextValue $ws.Range("D35") '0.2242'
$ws.Range("E35").Value = '  -3.79%  '

# Row 36
Set-TextValue $ws.Range("D36") '0.06440'
$ws.Range("E36").Value = '  -3.06%  '

# Row 37
Set-TextValue $ws.Range("D37") '0.02288'
$ws.Range("E37").Value = '  -1.53%  '

# Row 38
Set-TextValue $ws.Range("D38") '4.997'
$ws.Range("E38").Value = '  -3.66%  '

# Row 39
Set-TextValue $ws.Range("D39") '8.440'
$ws.Range("E39").Value = '  -3.95%  '

# Row 40
$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue $ws.Range("D40") '0.6165'
$ws.Range("E40").Value = '  -3.67%  '

# Row 41
$ws.Range("B41").Value = 'WEMIXTOKEN'
$ws.Range("C41").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue $ws.Range("D41") '1.445'
$ws.Range("E41").Value = '  +3.60%  '

# Row 42
Set-TextValue $ws.Range("D42") '11.11'
$ws.Range("E42").Value = '  -5.21%  '

# Row 43
Set-TextValue $ws.Range("D43") '1.171'
$ws.Range("E43").Value = '  +1.04%  '

# Row 44
$ws.Range("E44").Value = '  +0.15%  '

# Row 45
Set-TextValue $ws.Range("D45") '13.33'
$ws.Range("E45").Value = '  -2.36%  '

# Row 46
Set-TextValue $ws.Range("D46") '3.669'
$ws.Range("E46").Value = '  -0.28%  '

# Row 47
Set-TextValue $ws.Range("D47") '0.5755'
$ws.Range("E47").Value = '  -3.67%  '

# Row 48
Set-TextValue $ws.Range("D48") '124.69'
$ws.Range("E48").Value = '  -1.53%  '

# Row 49
$ws.Range("E49").Value = '  +2.66%  '

# Row 50
Set-TextValue $ws.Range("D50") '1.917'
$ws.Range("E50").Value = '  -4.13%  '

# Row 51
Set-TextValue $ws.Range("D51") '0.06816'
$ws.Range("E51").Value = '  -1.27%  '

